$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with freshly scraped values.
# Some Price values look numeric (e.g. "214.91"); force them to stay text
# (matching the original inline-string cell type) via a Text number format,
# then restore the cell style to Normal so no stray formatting is left behind.

$ws.Range("D2").Value = '26.991.23'
$ws.Range("E2").Value = '  -0.80%  '
$ws.Range("D3").Value = '1.620.02'
$ws.Range("E3").Value = '  -1.54%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.91'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.16%  '
$ws.Range("E6").Value = '  +0.25%  '
$ws.Range("E7").Value = '  +0.27%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.252'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.67%  '
$ws.Range("E9").Value = '  -0.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.12'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.97%  '
$ws.Range("E11").Value = '  -0.05%  '
$ws.Range("D12").Value = '1.624.48'
$ws.Range("E12").Value = '  -1.44%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.11'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.87%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.540'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.35%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '64.40'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.84%  '
$ws.Range("D16").Value = '26.984.29'
$ws.Range("E16").Value = '  -0.78%  '
$ws.Range("D17").Value = '0.0₃0739'
$ws.Range("E17").Value = '  +0.15%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '215.66'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.66%  '
$ws.Range("E19").Value = '  +0.18%  '
$ws.Range("E20").Value = '  +0.83%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.34'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.53%  '
$ws.Range("E22").Value = '  -5.77%  '
$ws.Range("E23").Value = '  -2.59%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '147.35'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.45%  '
$ws.Range("E25").Value = '  +0.36%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.26'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.117'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.95%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.52'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.47%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0502'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.06%  '
$ws.Range("E30").Value = '  -0.97%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.34'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.02%  '
$ws.Range("E32").Value = '  -1.88%  '
$ws.Range("D33").Value = '1.331.37'
$ws.Range("E33").Value = '  +5.50%  '
$ws.Range("E34").Value = '  -1.16%  '
$ws.Range("E35").Value = '  +0.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0175'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.28%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.542'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.844'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.70%  '
$ws.Range("E39").Value = '  +0.23%  '
$ws.Range("E40").Value = '  +0.13%  '
$ws.Range("E41").Value = '  -0.94%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '64.49'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.36%  '
$ws.Range("D44").Value = '1.757.92'
$ws.Range("E44").Value = '  -1.60%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '90.33'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.40%  '
$ws.Range("E46").Value = '  -0.29%  '
$ws.Range("D47").Value = '0.0₆0106'
$ws.Range("E47").Value = '  -1.10%  '
$ws.Range("E48").Value = '  +22.43%  '
$ws.Range("E49").Value = '  -0.18%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0981'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.88%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.54'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.59%  '
